$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pairwise_comp")

$newText = "Enter pairwise comparisons in the white cells of the table or numerical data in the green cells. For the Direct Values column, if the smallest value is best, invert the value before entering it (e.g., `$10 as =1/10) ."

$cells = @("A2", "A11", "A19", "A26", "A34", "A42", "A50")
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).Value = $newText
}
